$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 4230
$ws1.Range("F8").Value = 4230
$ws1.Range("F12").Value = 6152
$ws1.Range("F19").Value = 9270
$ws1.Range("F21").Value = 2499
$ws1.Range("F22").Value = 195
$ws1.Range("F23").Value = 2325
$ws1.Range("F24").Value = 2471
$ws1.Range("F25").Value = 1401
$ws1.Range("F26").Value = 245
$ws1.Range("F30").Value = 334
$ws1.Range("F36").Value = 384
$ws1.Range("F41").Value = 245

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 4230
$ws4.Range("F18").Value = 6152
$ws4.Range("F22").Value = 9270
$ws4.Range("F24").Value = 2499
$ws4.Range("F25").Value = 195
$ws4.Range("F26").Value = 2325
$ws4.Range("F27").Value = 2471
$ws4.Range("F28").Value = 1401
$ws4.Range("F29").Value = 245
$ws4.Range("F33").Value = 334
$ws4.Range("F37").Value = 384
$ws4.Range("F41").Value = 245
